# Updated partnership process:
#  - Fixed the bug in the partnership process U2b
#  - Partnership estimates were updated too
#
# Target sheet: ColumnsNumberParameters
#   Row 24 (columnsPartnershipU1b): B24 30 -> "31" (entered as text, quote-prefixed)
#   Row 25 (columnsPartnershipU2b): B25 33 -> "35" (entered as text, quote-prefixed)
#   Selection moves from B32 to B26

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ColumnsNumberParameters")

# Update the partnership estimate column numbers. A leading apostrophe enters the
# value with Excel's "quote prefix" so it is stored as text (matching the
# existing style used by the other text-typed column numbers on this sheet).
$ws.Range("B24").Value = "'31"
$ws.Range("B25").Value = "'35"

# Move the active selection to B26 (as last left by the author when saving).
$ws.Activate()
$ws.Range("B26").Select()
